$wb = $excel.ActiveWorkbook

# --- Sheet "LIST": mirror column A (rows 2-14) into column D, keeping the ---
# --- same text-format styling used on the source cells in column A.      ---
$wsList = $wb.Worksheets.Item("LIST")
for ($r = 2; $r -le 14; $r++) {
    $srcCell = $wsList.Cells.Item($r, 1)
    $dstCell = $wsList.Cells.Item($r, 4)
    $dstCell.Value = $srcCell.Value2
    if ($r -le 11) {
        $dstCell.NumberFormat = "@"
    }
}

# --- Sheet "Feuil1": move the selection to B21 ---
$wsFeuil1 = $wb.Worksheets.Item("Feuil1")
$wsFeuil1.Range("B21").Select()

# --- Back on "LIST": move the selection to C19 (keeps LIST the active tab) ---
$wsList.Range("C19").Select()
